$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, A(Sending cluster), B(Ligand symbol), C(Receptor symbol), D(Target cluster),
# E..T = Ligand-expressing cells ... Edge total expression derived specificity
$newRows = @(
    @(2, "ECs", "Sema3d", "Nrp1", "ECs", 1, 0.3333333333333333, 0.01606133333333333, 0.048184, 0.001454104962996282, 0.001454104962996282, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 0.8498100239582222, 7.648290215623999, 0.0006051134823402718, 0.0006051134823402717),
    @(3, "ECs", "Sema3d", "Nrp1", "FAPs", 1, 0.3333333333333333, 0.01606133333333333, 0.048184, 0.001454104962996282, 0.001454104962996282, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 0.7590595502444444, 6.831535952199999, 0.0005404939395897698, 0.0005404939395897698),
    @(4, "ECs", "Sema3d", "Nrp1", "sCs", 1, 0.3333333333333333, 0.01606133333333333, 0.048184, 0.001454104962996282, 0.001454104962996282, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 0.4332481599164445, 3.899233439248, 0.0003084975410662407, 0.0003084975410662407),
    @(5, "FAPs", "Sema3d", "Nrp1", "ECs", 3, 1, 8.111092999999999, 24.333279, 0.7343338402762993, 0.7343338402762993, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 429.1603936985742, 3862.443543287168, 0.3055868170439857, 0.3055868170439857),
    @(6, "FAPs", "Sema3d", "Nrp1", "FAPs", 3, 1, 8.111092999999999, 24.333279, 0.7343338402762993, 0.7343338402762993, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 383.3307283270915, 3449.976554943824, 0.2729534665002286, 0.2729534665002286),
    @(7, "FAPs", "Sema3d", "Nrp1", "sCs", 3, 1, 8.111092999999999, 24.333279, 0.7343338402762993, 0.7343338402762993, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 218.7935487191486, 1969.141938472338, 0.1557935567320852, 0.1557935567320852),
    @(8, "sCs", "Sema3d", "Nrp1", "ECs", 3, 1, 2.918357333333333, 8.755072, 0.2642120547607045, 0.2642120547607045, 3, 1, 52.91030366666666, 158.730911, 0.4161415425564564, 0.4161415425564564, 154.4111727145102, 1389.700554430592, 0.1099496120301305, 0.1099496120301305),
    @(9, "sCs", "Sema3d", "Nrp1", "FAPs", 3, 1, 2.918357333333333, 8.755072, 0.2642120547607045, 0.2642120547607045, 3, 1, 47.26005833333333, 141.780175, 0.3717021489810786, 0.3717021489810786, 137.9217378108444, 1241.2956402976, 0.09820818854126029, 0.09820818854126029),
    @(10, "sCs", "Sema3d", "Nrp1", "sCs", 3, 1, 2.918357333333333, 8.755072, 0.2642120547607045, 0.2642120547607045, 3, 1, 26.97460733333333, 80.923822, 0.2121563084624651, 0.2121563084624651, 78.72154312502045, 708.4938881251841, 0.05605425418931376, 0.05605425418931376)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $newRows) {
    $rowNum = $entry[0]
    for ($i = 1; $i -lt $entry.Count; $i++) {
        $colLetter = $cols[$i - 1]
        $ws.Range("$colLetter$rowNum").Value = $entry[$i]
    }
}
